$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.00000003236075771603351
$ws.Range("E2").Value = 0.00000003236075771603351

# Row 3
$ws.Range("D3").Value = 0.9999998467278153
$ws.Range("E3").Value = 0.9999998467278153

# Row 4
$ws.Range("D4").Value = 0.005966157138455723
$ws.Range("E4").Value = 0.005966157138455723

# Row 5
$ws.Range("D5").Value = 0.0001203705641323416
$ws.Range("E5").Value = 0.0001203705641323416

# Row 6
$ws.Range("D6").Value = 0.1211498380866866
$ws.Range("E6").Value = 0.1211498380866866

# Row 7
$ws.Range("D7").Value = 0.5507538472793829
$ws.Range("E7").Value = 0.4492461527206171

# Row 8
$ws.Range("D8").Value = 0.999999262505412
$ws.Range("E8").Value = 0.0000007374945879545791

# Row 9
$ws.Range("D9").Value = 0.9580982088655523
$ws.Range("E9").Value = 0.0419017911344477

# Row 10
$ws.Range("D10").Value = 0.9994948369111974
$ws.Range("E10").Value = 0.0005051630888025516

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.00000000003980581445873468
$ws.Range("E11").Value = 0.9999999999601942
$ws.Range("F11").Value = 4.041307926177979
$ws.Range("G11").Value = 0.8

# Row 12
$ws.Range("D12").Value = 0.0000000000003973903236973213
$ws.Range("E12").Value = 0.0000000000003973903236973213

# Row 13
$ws.Range("D13").Value = 0.9999999999965046
$ws.Range("E13").Value = 0.9999999999965046

# Row 14
$ws.Range("D14").Value = 0.0004844675398954387
$ws.Range("E14").Value = 0.0004844675398954387

# Row 15
$ws.Range("D15").Value = 0.00001038211105392742
$ws.Range("E15").Value = 0.00001038211105392742

# Row 16
$ws.Range("D16").Value = 0.08343829025094397
$ws.Range("E16").Value = 0.08343829025094397

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.4610201573135929
$ws.Range("E17").Value = 0.5389798426864071

# Row 18
$ws.Range("D18").Value = 0.999999995579286
$ws.Range("E18").Value = 0.000000004420714017072669

# Row 19
$ws.Range("D19").Value = 0.9969776422140343
$ws.Range("E19").Value = 0.003022357785965712

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.0000000000000001597667204035653
$ws.Range("E20").Value = 0.9999999999999999

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.0000000000000003920272621252376
$ws.Range("E21").Value = 0.9999999999999996
$ws.Range("F21").Value = 9.909256935119629
$ws.Range("G21").Value = 0.6
